$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I = "I0", J = "IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new
# header cells so they share the same bold / bordered / centered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I (I0) and J (IF) columns, rows 2-33
$data = @(
    @(8, 9),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(4, 5),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(5, 6),
    @(6, 6),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(4, 5),
    @(8, 8),
    @(5, 6),
    @(3, 3),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
